# Applies the cryptos.xlsx data refresh described by the commit diff.
# Prices (column D), volume deltas (column E), and for five rows the
# coin name/link (columns B/C) since those rows were re-ranked/swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'63.686.61"
$ws.Range("E2").Value = "  +3.27%  "
# Row 3
$ws.Range("D3").Value = "'2.945.05"
$ws.Range("E3").Value = "  +1.57%  "
# Row 4
$ws.Range("E4").Value = "  +0.12%  "
# Row 5
$ws.Range("D5").Value = "'573.91"
$ws.Range("E5").Value = "  -1.80%  "
# Row 6
$ws.Range("D6").Value = "'148.29"
$ws.Range("E6").Value = "  +0.83%  "
# Row 7
$ws.Range("E7").Value = "  -0.06%  "
# Row 8
$ws.Range("D8").Value = "'2.943.54"
$ws.Range("E8").Value = "  +1.60%  "
# Row 9
$ws.Range("E9").Value = "  -0.06%  "
# Row 10
$ws.Range("D10").Value = "'6.98"
$ws.Range("E10").Value = "  +4.13%  "
# Row 11
$ws.Range("E11").Value = "  +0.10%  "
# Row 12
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  +1.68%  "
# Row 13
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("E13").Value = "  -0.20%  "
# Row 14
$ws.Range("D14").Value = "'33.74"
$ws.Range("E14").Value = "  +2.92%  "
# Row 15
$ws.Range("E15").Value = "  +0.48%  "
# Row 16
$ws.Range("D16").Value = "'3.431.40"
$ws.Range("E16").Value = "  +1.49%  "
# Row 17
$ws.Range("D17").Value = "'63.624.40"
$ws.Range("E17").Value = "  +3.14%  "
# Row 18
$ws.Range("D18").Value = "'6.78"
$ws.Range("E18").Value = "  +2.05%  "
# Row 19
$ws.Range("D19").Value = "'2.937.95"
$ws.Range("E19").Value = "  +1.54%  "
# Row 20
$ws.Range("D20").Value = "'444.59"
$ws.Range("E20").Value = "  +2.03%  "
# Row 21
$ws.Range("D21").Value = "'13.47"
$ws.Range("E21").Value = "  +1.29%  "
# Row 22
$ws.Range("D22").Value = "'0.670"
$ws.Range("E22").Value = "  +1.88%  "
# Row 23
$ws.Range("D23").Value = "'7.05"
$ws.Range("E23").Value = "  +1.65%  "
# Row 24
$ws.Range("D24").Value = "'79.67"
$ws.Range("E24").Value = "  -0.30%  "
# Row 25
$ws.Range("D25").Value = "'12.14"
$ws.Range("E25").Value = "  +0.81%  "
# Row 26
$ws.Range("D26").Value = "'10.52"
$ws.Range("E26").Value = "  +3.09%  "
# Row 27
$ws.Range("D27").Value = "'2.15"
$ws.Range("E27").Value = "  +4.59%  "
# Row 28
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.03%  "
# Row 29
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "  +2.96%  "
# Row 30
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000106"
$ws.Range("E30").Value = "  -3.65%  "
# Row 31
$ws.Range("E31").Value = "  -0.69%  "
# Row 32
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -0.07%  "
# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.108"
$ws.Range("E33").Value = "  +0.47%  "
# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.28"
$ws.Range("E34").Value = "  +2.11%  "
# Row 35
$ws.Range("E35").Value = "  +0.09%  "
# Row 36
$ws.Range("D36").Value = "'0.964"
$ws.Range("E36").Value = "  +0.40%  "
# Row 37
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'2.11"
$ws.Range("E37").Value = "  +6.46%  "
# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'5.57"
$ws.Range("E38").Value = "  +1.58%  "
# Row 39
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  -2.65%  "
# Row 40
$ws.Range("D40").Value = "'49.12"
$ws.Range("E40").Value = "  +0.12%  "
# Row 41
$ws.Range("D41").Value = "'43.35"
$ws.Range("E41").Value = "  +13.98%  "
# Row 42
$ws.Range("E42").Value = "  +1.11%  "
# Row 43
$ws.Range("D43").Value = "'8.28"
$ws.Range("E43").Value = "  -0.53%  "
# Row 44
$ws.Range("D44").Value = "'0.282"
$ws.Range("E44").Value = "  +4.72%  "
# Row 45
$ws.Range("D45").Value = "'2.738.84"
$ws.Range("E45").Value = "  +2.09%  "
# Row 46
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'374.21"
$ws.Range("E46").Value = "  +9.27%  "
# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0346"
$ws.Range("E47").Value = "  +3.19%  "
# Row 48
$ws.Range("D48").Value = "'133.14"
$ws.Range("E48").Value = "  -1.30%  "
# Row 49
$ws.Range("E49").Value = "  -0.01%  "
# Row 50
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.104"
$ws.Range("E50").Value = "  +1.03%  "
# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000215"
$ws.Range("E51").Value = "  +7.92%  "
